# Lab 11 - Fix Problem 8 - While loop
# Turns the prior do-while-style code (label after the body, conditional
# branch back to the top) into a proper while loop: a start label before
# the condition check, a conditional branch to an end label, and an
# unconditional branch back to the start label, followed by the end label.

$d = $word.ActiveDocument

# --- 1. Insert ": startWhileLabel1" into the blank paragraph that sits
#        right before "load r1, a" (the start of the condition-check code).
#        Anchor on the (unique) ": loop1" label paragraph and walk backwards
#        to the "load r1, a" / "load r2, b" pair and the blank line above it,
#        since "load r1, a" alone appears several times earlier in the doc.
$count = $d.Paragraphs.Count
$loopLabelIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text.TrimEnd("`r")
    if ($ptext -eq ": loop1") {
        $loopLabelIndex = $i
        break
    }
}

$startPara = $null
if ($loopLabelIndex -gt 1) {
    for ($j = $loopLabelIndex - 1; $j -ge 1; $j--) {
        $ptext = $d.Paragraphs.Item($j).Range.Text.TrimEnd("`r")
        if ($ptext -eq "") {
            $startPara = $d.Paragraphs.Item($j)
            break
        }
        if ($ptext -ne "load r1, a" -and $ptext -ne "load r2, b") {
            break
        }
    }
}

if ($startPara -ne $null) {
    $r = $startPara.Range
    $r.Collapse(1)
    $r.InsertBefore(": startWhileLabel1") | Out-Null

    $fixRange = $startPara.Range
    $fixRange.Font.Name = "Times New Roman"
    $fixRange.Font.NameBi = "Times New Roman"
    $fixRange.Font.Size = 12
    $fixRange.Font.SizeBi = 12
}

# --- 2. ": loop1" becomes the conditional branch to the new end label.
$d.Content.Find.Execute(": loop1", $false, $false, $false, $false, $false, `
    $true, 1, $false, "branchEqual r1, r2, endWhileLabel1", 2) | Out-Null

# --- 3. "branchequal r1, r2, loop1" becomes the unconditional branch back
#        to the start label, and a new paragraph with ": endWhileLabel1"
#        is appended after it.
$d.Content.Find.Execute("branchequal r1, r2, loop1", $false, $false, $false, `
    $false, $false, $true, 1, $false, "branch startWhileLabel1", 2) | Out-Null

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter() | Out-Null
$newLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$newLast.Range.InsertBefore(": endWhileLabel1") | Out-Null
